$wb = $excel.ActiveWorkbook

# --- Update the summary text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 1.85 = 6699.86 pesos"), "✅ 1000 Bs = 1.85 = 6701.8 pesos"
$text = $text -replace [regex]::Escape("✅ 6699.86 pesos = 1.85 = 971.51 Bs"), "✅ 6701.8 pesos = 1.84 = 960.41 Bs"
$cell.Value2 = $text

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 3621.05
$wsTasas.Range("N12").Value = 3638.99
$wsTasas.Range("O12").Value = 521.4880000000001
